# edit.ps1 - apply 'literacy, 240930 modified 4' commit:
# Adds 10 new survey-response rows (482-491) to the Google-Forms-exported
# sheet, re-styles the former last row (481) as a normal banded row, gives
# the new last row (491) the special bottom-border treatment, and grows the
# Form_Responses1 table to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Formatting first: snapshot row 481's current ('last row') border
#    treatment onto the new last row (491) before row 481 itself gets
#    re-styled as an ordinary row.
# ------------------------------------------------------------------
$ws.Range("A481:L481").Copy()
$ws.Range("A491:L491").PasteSpecial(-4122)
$ws.Range("L481").Copy()
$ws.Range("M491").PasteSpecial(-4122)

# Row 481 is no longer the last row: give it the regular alternating-band style.
$ws.Range("A479:L479").Copy()
$ws.Range("A481:L481").PasteSpecial(-4122)
$ws.Range("N479").Copy()
$ws.Range("N481").PasteSpecial(-4122)

# New rows 482-490: copy formatting from the matching stable template row
# (alternating band colour x whether the response used column M or N).
# Rows that end in column N copy A:L then N separately (A:N in one shot
# leaves a stray, style-less placeholder cell in the always-skipped M
# column of both the template and the target row).
$ws.Range("A480:M480").Copy()
$ws.Range("A482:M482").PasteSpecial(-4122)

$ws.Range("A477:M477").Copy()
$ws.Range("A483:M483").PasteSpecial(-4122)

$ws.Range("A478:L478").Copy()
$ws.Range("A484:L484").PasteSpecial(-4122)
$ws.Range("N478").Copy()
$ws.Range("N484").PasteSpecial(-4122)

$ws.Range("A477:M477").Copy()
$ws.Range("A485:M485").PasteSpecial(-4122)

$ws.Range("A478:L478").Copy()
$ws.Range("A486:L486").PasteSpecial(-4122)
$ws.Range("N478").Copy()
$ws.Range("N486").PasteSpecial(-4122)

$ws.Range("A477:M477").Copy()
$ws.Range("A487:M487").PasteSpecial(-4122)

$ws.Range("A480:M480").Copy()
$ws.Range("A488:M488").PasteSpecial(-4122)

$ws.Range("A479:L479").Copy()
$ws.Range("A489:L489").PasteSpecial(-4122)
$ws.Range("N479").Copy()
$ws.Range("N489").PasteSpecial(-4122)

$ws.Range("A480:M480").Copy()
$ws.Range("A490:M490").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2) Cell values / formulas for the 10 new survey responses.
#    (Row 481's own values are untouched by this commit.)
# ------------------------------------------------------------------
# Row 482
$ws.Cells.Item(482,1).Value = 45568.431765636575
$ws.Cells.Item(482,2).Value = 'ytkay05@gmail.com'
$ws.Cells.Item(482,3).Value = '데이터사이언스학부'
$ws.Cells.Item(482,4).Value = 20243214
$ws.Cells.Item(482,5).Value = '김영민'
$ws.Cells.Item(482,6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(482,7).Value = 0.1
$ws.Cells.Item(482,8).Value = '6:4'
$ws.Cells.Item(482,9).Value = '20분의 1'
$ws.Cells.Item(482,10).Value = '20만호, 69만명'
$ws.Cells.Item(482,11).Value = '충청'
$ws.Cells.Item(482,12).Value = 'Red'
$ws.Cells.Item(482,13).Value = '모름/무응답'

# Row 483
$ws.Cells.Item(483,1).Value = 45568.45962278935
$ws.Cells.Item(483,2).Value = 'raon02271@naver.com'
$ws.Cells.Item(483,3).Value = '언어청각학부'
$ws.Cells.Item(483,4).Value = 20233954
$ws.Cells.Item(483,5).Value = '이채윤'
$ws.Cells.Item(483,6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(483,7).Value = 0.1
$ws.Cells.Item(483,8).Value = '6:4'
$ws.Cells.Item(483,9).Value = '10분의 1'
$ws.Cells.Item(483,10).Value = '20만호, 69만명'
$ws.Cells.Item(483,11).Value = '충청'
$ws.Cells.Item(483,12).Value = 'Red'
$ws.Cells.Item(483,13).Value = '근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다.'

# Row 484
$ws.Cells.Item(484,1).Value = 45568.57926799769
$ws.Cells.Item(484,2).Value = 'krdevmon@gmail.com'
$ws.Cells.Item(484,3).Value = '빅데이터'
$ws.Cells.Item(484,4).Value = 20205198
$ws.Cells.Item(484,5).Value = '안봉근'
$ws.Cells.Item(484,6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(484,7).Value = 0.1
$ws.Cells.Item(484,8).Value = '4:6'
$ws.Cells.Item(484,9).Value = '30분의 1'
$ws.Cells.Item(484,10).Value = '20만호, 69만명'
$ws.Cells.Item(484,11).Value = '충청'
$ws.Cells.Item(484,12).Value = 'Black'
$ws.Cells.Item(484,14).Value = '노동자가 과도한 연장근로를 받을 수 있어 반대한다.'

# Row 485
$ws.Cells.Item(485,1).Value = 45568.61332790509
$ws.Cells.Item(485,2).Value = 'hyj13223@naver.com'
$ws.Cells.Item(485,3).Value = '정치행정학과'
$ws.Cells.Item(485,4).Value = 20212432
$ws.Cells.Item(485,5).Value = '이현진'
$ws.Cells.Item(485,6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(485,7).Value = 0.1
$ws.Cells.Item(485,8).Value = '6:4'
$ws.Cells.Item(485,9).Value = '20분의 1'
$ws.Cells.Item(485,10).Value = '20만호, 69만명'
$ws.Cells.Item(485,11).Value = '충청'
$ws.Cells.Item(485,12).Value = 'Red'
$ws.Cells.Item(485,13).Value = '근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다.'

# Row 486
$ws.Cells.Item(486,1).Value = 45568.61960128472
$ws.Cells.Item(486,2).Value = 'ndd1016@gmail.com'
$ws.Cells.Item(486,3).Value = '환경생명공학과'
$ws.Cells.Item(486,4).Value = 20213737
$ws.Cells.Item(486,5).Value = '최형렬'
$ws.Cells.Item(486,6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(486,7).Value = 0.1
$ws.Cells.Item(486,8).Value = '6:4'
$ws.Cells.Item(486,9).Value = '20분의 1'
$ws.Cells.Item(486,10).Value = '20만호, 69만명'
$ws.Cells.Item(486,11).Value = '충청'
$ws.Cells.Item(486,12).Value = 'Black'
$ws.Cells.Item(486,14).Value = '찬성한다.'

# Row 487
$ws.Cells.Item(487,1).Value = 45568.623457442125
$ws.Cells.Item(487,2).Value = 'ran8410@naver.com'
$ws.Cells.Item(487,3).Value = '언어병리학과'
$ws.Cells.Item(487,4).Value = 20203918
$ws.Cells.Item(487,5).Value = '박경란'
$ws.Cells.Item(487,6).Value = '과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다.'
$ws.Cells.Item(487,7).Value = 0.3
$ws.Cells.Item(487,8).Value = '6:4'
$ws.Cells.Item(487,9).Value = '15분의 1'
$ws.Cells.Item(487,10).Value = '20만호, 69만명'
$ws.Cells.Item(487,11).Value = '경기'
$ws.Cells.Item(487,12).Value = 'Red'
$ws.Cells.Item(487,13).Value = '근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다.'

# Row 488
$ws.Cells.Item(488,1).Value = 45568.66208879629
$ws.Cells.Item(488,2).Value = 'joon020978@gmail.com'
$ws.Cells.Item(488,3).Value = '사회학과'
$ws.Cells.Item(488,4).Value = 20242230
$ws.Cells.Item(488,5).Value = '이준'
$ws.Cells.Item(488,6).Value = '‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다.'
$ws.Cells.Item(488,7).Value = 0.3
$ws.Cells.Item(488,8).Value = '7:3'
$ws.Cells.Item(488,9).Value = '20분의 1'
$ws.Cells.Item(488,10).Value = '44만호, 153만명'
$ws.Cells.Item(488,11).Value = '평안'
$ws.Cells.Item(488,12).Value = 'Red'
$ws.Cells.Item(488,13).Value = '근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다.'

# Row 489
$ws.Cells.Item(489,1).Value = 45568.66246140046
$ws.Cells.Item(489,2).Value = 'tqwquqqi@naver.com'
$ws.Cells.Item(489,3).Value = '심리학과'
$ws.Cells.Item(489,4).Value = 20222905
$ws.Cells.Item(489,5).Value = '강하늘'
$ws.Cells.Item(489,6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(489,7).Value = 0.1
$ws.Cells.Item(489,8).Value = '6:4'
$ws.Cells.Item(489,9).Value = '20분의 1'
$ws.Cells.Item(489,10).Value = '20만호, 69만명'
$ws.Cells.Item(489,11).Value = '충청'
$ws.Cells.Item(489,12).Value = 'Black'
$ws.Cells.Item(489,14).Value = '찬성한다.'

# Row 490
$ws.Cells.Item(490,1).Value = 45568.66503364583
$ws.Cells.Item(490,2).Value = 'yongwoo7701@gmail.com'
$ws.Cells.Item(490,3).Value = '체육학과'
$ws.Cells.Item(490,4).Value = 20244130
$ws.Cells.Item(490,5).Value = '유용우'
$ws.Cells.Item(490,6).Value = '등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다.'
$ws.Cells.Item(490,7).Value = 0.7
$ws.Cells.Item(490,8).Value = '7:3'
$ws.Cells.Item(490,9).Value = '10분의 1'
$ws.Cells.Item(490,10).Value = '44만호, 153만명'
$ws.Cells.Item(490,11).Value = '충청'
$ws.Cells.Item(490,12).Value = 'Red'
$ws.Cells.Item(490,13).Value = '근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다.'

# Row 491
$ws.Cells.Item(491,1).Value = 45568.72029888889
$ws.Cells.Item(491,2).Value = 'sky0219msh@naver.com'
$ws.Cells.Item(491,3).Value = '환경생명공학과'
$ws.Cells.Item(491,4).Value = 20201108
$ws.Cells.Item(491,5).Value = '최하늘'
$ws.Cells.Item(491,6).Value = '실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다.'
$ws.Cells.Item(491,7).Value = 0.1
$ws.Cells.Item(491,8).Value = '5:5'
$ws.Cells.Item(491,9).Value = '10분의 1'
$ws.Cells.Item(491,10).Value = '44만호, 153만명'
$ws.Cells.Item(491,11).Value = '전라'
$ws.Cells.Item(491,12).Value = 'Red'
$ws.Cells.Item(491,13).Value = '근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다.'

# ------------------------------------------------------------------
# 3) Grow the Form_Responses1 table so the new rows are part of it.
# ------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:N491"))
